# Applies the "Second commit for the week 6 assessment files" edits:
#  1. Append the Heroku deployed-link URL after "Deployed link: ".
#  2. Merge the " (" / "locally: " runs in the first "(locally: ...)" note.
#  3/4. Merge the split "... This button should turn orange when it is
#       hovered over." sentences back into single runs.
#  5/6. Merge the split "Losses" sentences back into single runs.
#  7. Mark the screenshot picture's run as NoProofing (<w:noProof/>).

$d = $word.ActiveDocument

function Replace-InParagraph($doc, $paraIndex, $searchText, $replaceText) {
    # Find + in-place Replace: this is how Word naturally collapses a
    # phrase that currently spans several identically-formatted runs
    # back down into a single run.
    $p = $doc.Paragraphs($paraIndex)
    $pr = $p.Range
    $rng = $doc.Range($pr.Start, $pr.End)
    $ok = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output ("NOT FOUND (para " + $paraIndex + "): " + $searchText)
    }
}

# 1. "Deployed link: " -> "Deployed link: https://week6-assessment.herokuapp.com/"
$depPara = $d.Paragraphs(8)
$depRange = $d.Range($depPara.Range.Start, $depPara.Range.End - 1)
$depRange.Collapse(0)
$depRange.InsertAfter("https://week6-assessment.herokuapp.com/")

# 2. " (" + "locally: " -> " (locally: " (first occurrence, the one using rId6)
Replace-InParagraph $d 76 " (locally: " " (locally: "

# 3. " selection. " + "This button should turn " + "orange " + "when it is hovered over."
Replace-InParagraph $d 78 " selection. This button should turn orange when it is hovered over." " selection. This button should turn orange when it is hovered over."

# 4. ". " + "This button should turn orange when it is hovered over."
Replace-InParagraph $d 79 ". This button should turn orange when it is hovered over." ". This button should turn orange when it is hovered over."

# 5. "When a Player wins, the "Losses" are " + "updated incrementally" + " added."
Replace-InParagraph $d 122 "When a Player wins, the “Losses” are updated incrementally added." "When a Player wins, the “Losses” are updated incrementally added."

# 6. "When a Player losses, the "Losses" are properly " + "updated incrementally" + "."
Replace-InParagraph $d 123 "When a Player losses, the “Losses” are properly updated incrementally." "When a Player losses, the “Losses” are properly updated incrementally."

# 7. Add <w:noProof/> to the run that hosts the inline screenshot picture.
$shape = $d.InlineShapes.Item(1)
$shape.Range.NoProofing = $true
